# fix index conditioning issue
#
# Adds three new cpars rows (I_beta / VI_beta / SpI_beta) describing the
# beta parameters used to condition index observation error for
# Data@Ind, Data@VInd and Data@SpInd, inserted into the "Obs" group of
# the cpars sheet (right after the existing "hbias" row). Also nudges a
# couple of row heights on the Fleet sheet and leaves the workbook's
# view state pointing at the newly-edited area.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. cpars sheet: insert the three new rows under the "Obs" section,
#    immediately after the "hbias" row (row 76), pushing the
#    "BioEco" block and the trailing "Data" row down accordingly.
# ---------------------------------------------------------------
$cpars = $wb.Worksheets.Item("cpars")

$cpars.Rows("76:78").Insert()

# Column A (the three new slot names) first, then column C (their
# descriptions) - this is the order the shared-string table picks the
# new entries up in, matching how they'd land if typed in column order.
$cpars.Cells.Item(76, 1).Value = 'I_beta'
$cpars.Cells.Item(77, 1).Value = 'VI_beta'
$cpars.Cells.Item(78, 1).Value = 'SpI_beta'

$cpars.Cells.Item(76, 3).Value = 'Beta for hyperstability/depletion for `Data@Ind`'
$cpars.Cells.Item(77, 3).Value = 'Beta for hyperstability/depletion for `Data@VInd`'
$cpars.Cells.Item(78, 3).Value = 'Beta for hyperstability/depletion for `Data@SpInd`'

$cpars.Cells.Item(76, 2).Value = 'numeric vector length nsim'
$cpars.Cells.Item(77, 2).Value = 'numeric vector length nsim'
$cpars.Cells.Item(78, 2).Value = 'numeric vector length nsim'

$cpars.Cells.Item(76, 4).Value = 'Obs'
$cpars.Cells.Item(77, 4).Value = 'Obs'
$cpars.Cells.Item(78, 4).Value = 'Obs'

# ---------------------------------------------------------------
# 2. Fleet sheet: a couple of description rows grew taller (wrapped
#    text reflowed), bump their row heights to match.
# ---------------------------------------------------------------
$fleet = $wb.Worksheets.Item("Fleet")
$fleet.Rows.Item(19).RowHeight = 60
$fleet.Rows.Item(20).RowHeight = 45

# ---------------------------------------------------------------
# 3. View-state bookkeeping: leave the cursor parked on the Obs
#    sheet at D12, then finish with the cpars sheet active and
#    scrolled/selected near the newly-added rows - matching where
#    the author ended up after making the edit.
# ---------------------------------------------------------------
$obs = $wb.Worksheets.Item("Obs")
$obs.Activate()
$obs.Range("D12").Select()

$cpars.Activate()
$cpars.Range("C78").Select()
